# Corrected APDU responses handling
#
# This script edits Slide 2 of the BLE sequence-diagram deck:
#   * Renames the plain "OK" acknowledgement labels to "ACK".
#   * Fixes the "Notification: APDU Response(s) Ready" typo/size.
#   * Inserts a second "Read: APDU Responses" request (with its arrow)
#     between the "APDU Responses (fragment 0)" and "(fragment 1)" rows,
#     by duplicating the existing "Read: APDU Responses" request/arrow
#     and sliding it down.
#   * Removes the four stray client "OK" reply rows that are no longer
#     needed now that the read is modelled explicitly.
#   * Moves the remaining rows below the insertion point down to make
#     room, and re-styles/repositions the "APDU Responses (fragment 1)"
#     acknowledgement arrow.

$EMU = 12700.0  # EMU per point

function ToPt($emu) { return $emu / $EMU }

function Get-ShapeById {
    param($slide, [int]$id)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) { return $sh }
    }
    return $null
}

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

# ---------------------------------------------------------------------
# 1. Simple "OK" -> "ACK" text fixes (shapes that persist unchanged
#    otherwise).
# ---------------------------------------------------------------------
foreach ($id in 47, 69, 71, 75, 79, 96) {
    $sh = Get-ShapeById $s $id
    $sh.TextFrame.TextRange.Text = "ACK"
}

# ---------------------------------------------------------------------
# 2. "Notification: APDU Response Ready" -> "...Responses Ready", plus
#    the textbox shrinks back down to the single-line height used
#    everywhere else.
# ---------------------------------------------------------------------
$sh59 = Get-ShapeById $s 59
$sh59.TextFrame.TextRange.Text = "Notification: APDU Responses Ready"
$sh59.Left = ToPt 1642783
$sh59.Top = ToPt 4149007
$sh59.Width = ToPt 2851107
$sh59.Height = ToPt 246221

# ---------------------------------------------------------------------
# 3. Duplicate the existing "Read: APDU Responses" request (shape 61)
#    and its arrow (shape 60) to create a second read, inserted between
#    the fragment-0 and fragment-1 response rows.
# ---------------------------------------------------------------------
$connSrc = Get-ShapeById $s 60
$tbSrc = Get-ShapeById $s 61

$connNew = $connSrc.Duplicate()
$connNew.Name = "Straight Arrow Connector 49"
$connNew.Left = ToPt 615004
$connNew.Top = ToPt 5403554

$tbNew = $tbSrc.Duplicate()
$tbNew.Name = "TextBox 50"
$tbNew.Left = ToPt 1067556
$tbNew.Top = ToPt 5190533

# ---------------------------------------------------------------------
# 4. Move the connector/textbox pairs that sit at/after the original
#    "Read: APDU Responses" row down to make room for the inserted row
#    and to re-flow the rest of the diagram.
# ---------------------------------------------------------------------
$sh58 = Get-ShapeById $s 58
$sh58.Left = ToPt 586579
$sh58.Top = ToPt 4362028

$sh60 = Get-ShapeById $s 60
$sh60.Left = ToPt 586579
$sh60.Top = ToPt 4919276

$sh61 = Get-ShapeById $s 61
$sh61.Left = ToPt 1039131
$sh61.Top = ToPt 4706255

$sh62 = Get-ShapeById $s 62
$sh62.Left = ToPt 596578
$sh62.Top = ToPt 5159395

$sh63 = Get-ShapeById $s 63
$sh63.Left = ToPt 1642783
$sh63.Top = ToPt 4946374

$sh66 = Get-ShapeById $s 66
$sh66.Left = ToPt 615004
$sh66.Top = ToPt 6184173

$sh67 = Get-ShapeById $s 67
$sh67.Left = ToPt 1067556
$sh67.Top = ToPt 5971152

# ---------------------------------------------------------------------
# 5. Remove the four stray client "OK" reply rows (arrow + textbox each)
#    that used to follow the Notification / Read / fragment-0 /
#    fragment-1 rows.
# ---------------------------------------------------------------------
foreach ($id in 80, 81, 83, 84, 86, 87, 90, 91) {
    $sh = Get-ShapeById $s $id
    $sh.Delete()
}

# ---------------------------------------------------------------------
# 6. Reposition/restyle the "APDU Responses (fragment 1)" ack arrow and
#    its textbox, and move the trailing rows up to close the gap left
#    by the removed "OK" rows.
# ---------------------------------------------------------------------
$sh88 = Get-ShapeById $s 88
$sh88.Left = ToPt 613358
$sh88.Top = ToPt 5653603

$sh89 = Get-ShapeById $s 89
$sh89.Left = ToPt 1659563
$sh89.Top = ToPt 5440582

$sh95 = Get-ShapeById $s 95
$sh95.Left = ToPt 598139
$sh95.Top = ToPt 6376355

$sh96 = Get-ShapeById $s 96
$sh96.Left = ToPt 3649736
$sh96.Top = ToPt 6163334
